$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(1, "ARROZ PAISA SUBA", 46015, 166000),
  @(2, "CAMPO VERDE TOCANCIPA", 46021, 475000),
  @(3, "CAMPO VERDE ZIPAQUIRA", 46014, 71900),
  @(4, "CANTON WOK", 46015, 252000),
  @(5, "CARNES JOHANA", 46011, 166000),
  @(6, "CARNES JOHANA", 46021, 354000),
  @(7, "CIMARRON DORADO", 46010, 375000),
  @(8, "CIMARRON DORADO", 46020, 449800),
  @(9, "CLIENTE PAOLA", 46018, 174000),
  @(10, "CRISTIAN ACACIAS", 46009, 1000000),
  @(11, "DARWIN FUTBOL", 45921, 200000),
  @(12, "DAVIDCITO", 45947, 100000),
  @(13, "FRANCO", 46017, 545800),
  @(14, "FRANCO", 45996, 20000),
  @(15, "LA 13", 46021, 1137000),
  @(16, "LA CABAÑA", 46020, 215300),
  @(17, "LA PAMPA", 46006, 229900),
  @(18, "LA PROMESA", 46020, 151000),
  @(19, "LA SELECTA", 45912, 82000),
  @(20, "MAFE", 46017, 190000),
  @(21, "MERKA FRUVER ALEJANDRO", 46021, 1257600),
  @(22, "MERKA FRUVER DEXI", 45995, 339000),
  @(23, "MERKA FRUVER DEXI", 45988, 15400),
  @(24, "MICHAEL", 46011, 80000),
  @(25, "NEVADA", 46020, 195000),
  @(26, "PARAÍSO FUNZA", 46020, 276000),
  @(27, "PARAÍSO MOSQUERA", 46013, 328800),
  @(28, "PINILLA", 45931, 82000),
  @(29, "PLANADAS NUEVO", 46020, 88400),
  @(30, "PLAZA JESSICA", 46014, 1655400),
  @(31, "PUNTA DE ANCA", 46017, 507000),
  @(32, "SAMY 2", 46020, 304000),
  @(33, "SAMY 2", 46021, 203000),
  @(34, "SAMY 2", 46013, 142000),
  @(35, "SAN JOAQUIN", 46015, 229300),
  @(36, "SANTANDER SUR", 46014, 253000),
  @(37, "SANTANDER SUR", 46018, 218000),
  @(38, "TIMO", 46015, 189000),
  @(39, "WILINTONG", 46006, 150000)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 3).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $false
    $r = $r + 1
}

Write-Host "Rows written:" ($r - 2)